{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The second paragraph is the numbered item describing the protoboard / pin\n// connection; the new analysis paragraph is inserted right after it, as a\n// new sibling list item (same list numbering/style/formatting).\nconst target = paragraphs.items[1];\n\nconst newParagraph = target.insertParagraph(\n  \"La funci\u00f3n empezar\u00e1 iterando desde la posici\u00f3n de memoria del primer elemento de la matriz y se mover\u00e1 un bloque de memoria del tipo de dato usado por cada iteraci\u00f3n en la cual encender\u00e1 y apagar\u00e1 los leds seg\u00fan un tiempo de espera de 500ms y luego verificamos mediante las se\u00f1ales anal\u00f3gicas si el led si se enciende o no\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The second paragraph (numbered item about the \"pin digital\" / protoboard\n# connection) gets a new sibling list item right after it, with the same\n# list-paragraph formatting (style, numbering, font), containing the new\n# analysis text about how the function iterates memory.\n$target = $d.Paragraphs.Item(2)\n\n$rng = $target.Range\n$rng.Collapse(0)  # wdCollapseEnd\n$rng.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item(3)\n$newPara.Range.Text = \"La funci\u00f3n empezar\u00e1 iterando desde la posici\u00f3n de memoria del primer elemento de la matriz y se mover\u00e1 un bloque de memoria del tipo de dato usado por cada iteraci\u00f3n en la cual encender\u00e1 y apagar\u00e1 los leds seg\u00fan un tiempo de espera de 500ms y luego verificamos mediante las se\u00f1ales anal\u00f3gicas si el led si se enciende o no\"\n"}
